$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report row is inserted at row 141 ("Fruta / hortaliza, semanal"
# weekly update), pushing the existing rows 141-216 down to 142-217 and
# extending the sheet's used range from A1:T216 to A1:T217.
$ws.Range("A141").EntireRow().Insert()

# Columns A, B, C, E, F, G, H, I, J hold the same constant values on every
# data row of this sheet (market/product identifiers) - replicate them here.
$ws.Range("A141").Value() = 4
$ws.Range("B141").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value() = "Los Lagos"
$ws.Range("E141").Value() = 10
$ws.Range("F141").Value() = "Fruta"
$ws.Range("G141").Value() = 100103
$ws.Range("H141").Value() = "Frutos de hueso (carozo)"
$ws.Range("I141").Value() = 100103004
$ws.Range("J141").Value() = "Durazno"

# New row-specific data for the inserted record.
$ws.Range("D141").Value() = 44904
$ws.Range("K141").Value() = "Early Majestic"
$ws.Range("L141").Value() = "Primera"
$ws.Range("M141").Value() = 600
$ws.Range("N141").Value() = 18000
$ws.Range("O141").Value() = 19000
$ws.Range("P141").Value() = 18500
$ws.Range("Q141").Value() = "$/caja 14 kilos empedrada"
$ws.Range("R141").Value() = "Región de O'Higgins"
$ws.Range("S141").Value() = 1321
$ws.Range("T141").Value() = 14
